$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data (row 7), mirroring the existing student rows
$ws.Range("A7").Value = "student6"
$ws.Range("B7").Value = 1234
$ws.Range("C7").Value = "student"
$ws.Range("D7").Value = "CD"

# Move the active selection to E7, as recorded after entering the row's data
$ws.Range("E7").Select()
